$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2023-05-31 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-06-01 Thursday", 2) | Out-Null

# Update each arithmetic-problem cell in the table (row-major order matches source order)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "76+10="
$t.Cell(1, 2).Range.Text = "72-35="
$t.Cell(1, 3).Range.Text = "82-16="
$t.Cell(1, 4).Range.Text = "55+30="
$t.Cell(1, 5).Range.Text = "16+26="

$t.Cell(2, 1).Range.Text = "89-53="
$t.Cell(2, 2).Range.Text = "34-27="
$t.Cell(2, 3).Range.Text = "42+34="
$t.Cell(2, 4).Range.Text = "91+2="
$t.Cell(2, 5).Range.Text = "15+67="

$t.Cell(3, 1).Range.Text = "54-17="
$t.Cell(3, 2).Range.Text = "88+8="
$t.Cell(3, 3).Range.Text = "5-0="
$t.Cell(3, 4).Range.Text = "73-33="
$t.Cell(3, 5).Range.Text = "66+18="

$t.Cell(4, 1).Range.Text = "49-31="
$t.Cell(4, 2).Range.Text = "26+57="
$t.Cell(4, 3).Range.Text = "81-37="
$t.Cell(4, 4).Range.Text = "69+23="
$t.Cell(4, 5).Range.Text = "49-13="

$t.Cell(5, 1).Range.Text = "20+54="
$t.Cell(5, 2).Range.Text = "46-31="
$t.Cell(5, 3).Range.Text = "31+48="
$t.Cell(5, 4).Range.Text = "36+4="
$t.Cell(5, 5).Range.Text = "51+16="

$t.Cell(6, 1).Range.Text = "40+52="
$t.Cell(6, 2).Range.Text = "73-37="
$t.Cell(6, 3).Range.Text = "4+58="
$t.Cell(6, 4).Range.Text = "82-61="
$t.Cell(6, 5).Range.Text = "0+45="

$t.Cell(7, 1).Range.Text = "0+65="
$t.Cell(7, 2).Range.Text = "56+7="
$t.Cell(7, 3).Range.Text = "41+26="
$t.Cell(7, 4).Range.Text = "6+83="
$t.Cell(7, 5).Range.Text = "5+48="

$t.Cell(8, 1).Range.Text = "72+17="
$t.Cell(8, 2).Range.Text = "1+52="
$t.Cell(8, 3).Range.Text = "34+30="
$t.Cell(8, 4).Range.Text = "69-6="
$t.Cell(8, 5).Range.Text = "87-43="

$t.Cell(9, 1).Range.Text = "5+26="
$t.Cell(9, 2).Range.Text = "75-46="
$t.Cell(9, 3).Range.Text = "17+78="
$t.Cell(9, 4).Range.Text = "76-0="
$t.Cell(9, 5).Range.Text = "24+68="

$t.Cell(10, 1).Range.Text = "17+69="
$t.Cell(10, 2).Range.Text = "6+19="
$t.Cell(10, 3).Range.Text = "10+70="
$t.Cell(10, 4).Range.Text = "43-21="
$t.Cell(10, 5).Range.Text = "83-75="

$t.Cell(11, 1).Range.Text = "13+36="
$t.Cell(11, 2).Range.Text = "67-9="
$t.Cell(11, 3).Range.Text = "70+6="
$t.Cell(11, 4).Range.Text = "27-3="
$t.Cell(11, 5).Range.Text = "36+9="

$t.Cell(12, 1).Range.Text = "76-30="
$t.Cell(12, 2).Range.Text = "63+22="
$t.Cell(12, 3).Range.Text = "1+12="
$t.Cell(12, 4).Range.Text = "86-8="
$t.Cell(12, 5).Range.Text = "84-25="

$t.Cell(13, 1).Range.Text = "69-33="
$t.Cell(13, 2).Range.Text = "70+18="
$t.Cell(13, 3).Range.Text = "3+8="
$t.Cell(13, 4).Range.Text = "37+46="
$t.Cell(13, 5).Range.Text = "0+63="

$t.Cell(14, 1).Range.Text = "62-42="
$t.Cell(14, 2).Range.Text = "92-32="
$t.Cell(14, 3).Range.Text = "84-40="
$t.Cell(14, 4).Range.Text = "73+23="
$t.Cell(14, 5).Range.Text = "67-21="

$t.Cell(15, 1).Range.Text = "53-52="
$t.Cell(15, 2).Range.Text = "34-16="
$t.Cell(15, 3).Range.Text = "27+50="
$t.Cell(15, 4).Range.Text = "85-1="
$t.Cell(15, 5).Range.Text = "14-13="

$t.Cell(16, 1).Range.Text = "25+0="
$t.Cell(16, 2).Range.Text = "77-7="
$t.Cell(16, 3).Range.Text = "87-58="
$t.Cell(16, 4).Range.Text = "47-27="
$t.Cell(16, 5).Range.Text = "76-1="

$t.Cell(17, 1).Range.Text = "52-43="
$t.Cell(17, 2).Range.Text = "83+8="
$t.Cell(17, 3).Range.Text = "59-32="
$t.Cell(17, 4).Range.Text = "83-29="
$t.Cell(17, 5).Range.Text = "91-18="

$t.Cell(18, 1).Range.Text = "10+68="
$t.Cell(18, 2).Range.Text = "7+12="
$t.Cell(18, 3).Range.Text = "5+4="
$t.Cell(18, 4).Range.Text = "84-49="
$t.Cell(18, 5).Range.Text = "69-57="

$t.Cell(19, 1).Range.Text = "97+0="
$t.Cell(19, 2).Range.Text = "12+12="
$t.Cell(19, 3).Range.Text = "45-34="
$t.Cell(19, 4).Range.Text = "34+36="
$t.Cell(19, 5).Range.Text = "19-17="

$t.Cell(20, 1).Range.Text = "33+9="
$t.Cell(20, 2).Range.Text = "84-16="
$t.Cell(20, 3).Range.Text = "76-5="
$t.Cell(20, 4).Range.Text = "11-8="
$t.Cell(20, 5).Range.Text = "14+54="
